$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Fix up header text and fill in the previously-missing G/I (New Cases /
#    New Deaths) values for the existing rows, then append the new row 9.
# ---------------------------------------------------------------------------
$ws.Range("L1").Value = "MasterSheet RowNo."

# TotalConfirmedNewCases (G) and TotalNewDeaths (I) for rows 2-8
$ws.Range("G2").Value = 2
$ws.Range("I2").Value = 0

$ws.Range("G3").Value = 1
$ws.Range("I3").Value = 0

$ws.Range("G4").Value = 0
$ws.Range("I4").Value = 0

$ws.Range("G5").Value = 2
$ws.Range("I5").Value = 0

$ws.Range("G6").Value = 0
$ws.Range("I6").Value = 0

$ws.Range("G7").Value = 3
$ws.Range("I7").Value = 0

$ws.Range("G8").Value = 0
$ws.Range("I8").Value = 0

# New row 9
$ws.Range("A9").Value = 71
$ws.Range("B9").Value = 247
$ws.Range("C9").Value = "SOUTHEAST ASIAN"
$ws.Range("D9").Value = 43921
$ws.Range("E9").Value = "Myanmar"
$ws.Range("F9").Value = 10
$ws.Range("G9").Value = 2
$ws.Range("H9").Value = 0
$ws.Range("I9").Value = 0
$ws.Range("J9").Value = "Local transmission"
$ws.Range("K9").Value = 0
$ws.Range("L9").Value = 5310

# ---------------------------------------------------------------------------
# 2. Column widths: A..O all 27 characters wide.
# ---------------------------------------------------------------------------
$ws.Range("A1:O1").ColumnWidth = 26.166666666666668

# ---------------------------------------------------------------------------
# 3. Formatting: center every used cell (A1:O9) horizontally & vertically,
#    and give the Date column (D) its own date number format.
#    A helper cell well outside the used range is used to build the combined
#    alignment format in a single step before being pasted onto the whole
#    range; this keeps the produced style table tight (matching how Excel
#    itself would coalesce identical formats into one style record).
# ---------------------------------------------------------------------------
$tmpl = $ws.Range("Z1")
$tmpl.HorizontalAlignment = -4108
$tmpl.VerticalAlignment = -4108

$tmpl.Copy()
$ws.Range("A1:O9").PasteSpecial(-4122)

$ws.Range("D1:D9").NumberFormat = "yyyy-mm-dd;"

$ws.Range("Z1").Delete(-4159)

$excel.CutCopyMode = 0
